$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-08-27 20:58:02"

for ($row = 2; $row -le 72; $row++) {
    $ws.Range("O$row").Value = $newTimestamp
}
